$d = $word.ActiveDocument

function ReplaceIn($range, $findText, $replaceText) {
    $range.Find.ClearFormatting()
    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. Title: "CS411 A2 Project Pitch" -> "CS411 A2 Project Proposal"
ReplaceIn $d.Paragraphs(1).Range "CS411 A2 Project Pitch" "CS411 A2 Project Proposal"

# 2. Authors paragraph: drop the spell-check run splits, "and" -> "&"
ReplaceIn $d.Paragraphs(2).Range "Daniel Mboweni, Sabina Razak, Deijah Lee-Carroll, and Kelly Zhang" "Daniel Mboweni, Sabina Razak, Deijah Lee-Carroll & Kelly Zhang"

# 2b. Move the (collapsed) "_GoBack" bookmark to the end of the authors paragraph.
# Bookmarks.Add on a truly collapsed (zero-length) range is unreliable in this host, so
# temporarily insert a placeholder character, wrap it with the bookmark, then delete the
# placeholder -- the bookmark collapses in place, anchored right after "Kelly Zhang".
$p2End = $d.Paragraphs(2).Range.End - 1
$d.Range($p2End, $p2End).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($p2End, $p2End + 1)) | Out-Null
$d.Range($p2End, $p2End + 1).Text = ""

# 3. Empty paragraph (between authors and "Idea 1:"): bold the paragraph mark
$d.Paragraphs(3).Range.Select()
$word.Selection.Bold = 1

# 4. "Idea #1:" -> "Idea 1:"
ReplaceIn $d.Paragraphs(4).Range "Idea #1:" "Idea 1:"

# 5. Idea 1 body: wording rewrite
$idea1Old = "We will be designing a web application that helps you locate nearby people with whom to play a sport, such as catch, basketball, soccer. etc. Based on the user’s location, activity preferences, skill level, and schedule availability, the system would send notifications of nearby pick-up games. The user would create a profile and log on using their Facebook or Gmail as a form of third party authentication. Profile information and preferences would be stored in a database that we create. Users will be able to communicate through a messaging platform, created with the Slack API, and the formed group can choose the venue of where to meet up using the Yelp API to check for reviews. The goal of this web application is to gather people for a spontaneous game of pick up."
$idea1New = "We will be designing a web application that helps you locate nearby people with whom to play a sport, such as catch, basketball, soccer, etc. Based on the user’s location, activity preferences, skills, and schedule availability, the app would notify interested users in nearby pick-up games. The user would create a profile and log on using their Facebook or Gmail as a form of third party authentication. Their profile information and selected preferences would be stored in a database that we create. Users will be able to communicate through a message platform based on the Slack API. The formed group can choose the venue of where to meet up using the Yelp API to check for reviews. The goal of this web app is to gather people for a spontaneous game of pick-up."
ReplaceIn $d.Paragraphs(5).Range $idea1Old $idea1New

# 6. "Idea #2:" -> "Idea 2:"
ReplaceIn $d.Paragraphs(7).Range "Idea #2:" "Idea 2:"

# 7. Idea 2 body: wording rewrite (this also removes the original "_GoBack" bookmark
#    that used to live in this paragraph, since that text run is fully replaced)
$idea2Old = "We will be designing a web application that builds a playlist of songs based on your current mood, location, and/or weather. The user logs on using Facebook or Gmail as a third-party authentication. The database that is constructed would store the user’s profile information, as well as past playlists put together, with the name of the songs chosen for that specific day to create a history. Once logged in, user enters their mood, and that accompanied with the weather and location (if accessed) would utilize Spotify to build the playlist for that user. The APIs we would use for this are Yahoo! Weather, Spotify, and Google Cloud Natural (for sentiment analysis). Google Cloud Natural would be used to analyze text messages and predict the user’s mood if desired."
$idea2New = "We will be designing a web application that builds a playlist of songs based on your current mood, location, and/or weather. The user logs on using Facebook, or Gmail as a third-party authentication. The constructed database would store the user’s profile information, as well as past playlists put together to form a history, with the name of the songs chosen for that specific day. Once logged in, users enter their mood, and that accompanied with the weather and location (if accessed) would utilize Spotify to build a playlist. The APIs we would use for this are Yahoo! Weather, Spotify, and Google Cloud Natural. Google Cloud Natural would be used to analyze text messages and predict the user’s mood, if desired, through sentiment analysis."
ReplaceIn $d.Paragraphs(8).Range $idea2Old $idea2New

Write-Host "Edits applied"
